$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrBF = New-Object "object[,]" 24,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.0374894944573
$arrBF[0,2] = 1.03936817208185
$arrBF[0,3] = 1.041140036502686
$arrBF[0,4] = 1.046289641449599
$arrBF[1,0] = 1.019999999999999
$arrBF[1,1] = 1.038566811896721
$arrBF[1,2] = 1.040169012468732
$arrBF[1,3] = 1.042166034829581
$arrBF[1,4] = 1.047584052641385
$arrBF[2,0] = 1.02
$arrBF[2,1] = 1.03926383179843
$arrBF[2,2] = 1.040687117658769
$arrBF[2,3] = 1.042830196718068
$arrBF[2,4] = 1.048421906658479
$arrBF[3,0] = 1.02
$arrBF[3,1] = 1.039556841426265
$arrBF[3,2] = 1.040904906915272
$arrBF[3,3] = 1.043109475857804
$arrBF[3,4] = 1.048774210006001
$arrBF[4,0] = 1.02
$arrBF[4,1] = 1.039606038026338
$arrBF[4,2] = 1.040941473386498
$arrBF[4,3] = 1.043156371902219
$arrBF[4,4] = 1.048833367413422
$arrBF[5,0] = 1.02
$arrBF[5,1] = 1.039267747075345
$arrBF[5,2] = 1.040690027856718
$arrBF[5,3] = 1.04283392820195
$arrBF[5,4] = 1.048426613877792
$arrBF[6,0] = 1.02
$arrBF[6,1] = 1.03785359524314
$arrBF[6,2] = 1.039638838837819
$arrBF[6,3] = 1.04148672121338
$arrBF[6,4] = 1.046727036090628
$arrBF[7,0] = 1.02
$arrBF[7,1] = 1.035361057715142
$arrBF[7,2] = 1.037785808838976
$arrBF[7,3] = 1.039114843417998
$arrBF[7,4] = 1.043734258787805
$arrBF[8,0] = 1.02
$arrBF[8,1] = 1.03369889634936
$arrBF[8,2] = 1.036549988055411
$arrBF[8,3] = 1.037534964314455
$arrBF[8,4] = 1.041740376298539
$arrBF[9,0] = 1.02
$arrBF[9,1] = 1.032979038334871
$arrBF[9,2] = 1.036014753035519
$arrBF[9,3] = 1.036851179214901
$arrBF[9,4] = 1.040877287593327
$arrBF[10,0] = 1.02
$arrBF[10,1] = 1.032711630138069
$arrBF[10,2] = 1.035815925309119
$arrBF[10,3] = 1.036597237416965
$arrBF[10,4] = 1.040556737217633
$arrBF[11,0] = 1.02
$arrBF[11,1] = 1.032768991067629
$arrBF[11,2] = 1.035858575330019
$arrBF[11,3] = 1.036651706702892
$arrBF[11,4] = 1.040625494622923
$arrBF[12,0] = 1.02
$arrBF[12,1] = 1.032956934722592
$arrBF[12,2] = 1.035998318233495
$arrBF[12,3] = 1.036830187362839
$arrBF[12,4] = 1.04085079001238
$arrBF[13,0] = 1.02
$arrBF[13,1] = 1.033072730209708
$arrBF[13,2] = 1.036084416111733
$arrBF[13,3] = 1.036940161323331
$arrBF[13,4] = 1.040989607065789
$arrBF[14,0] = 1.02
$arrBF[14,1] = 1.033746668091335
$arrBF[14,2] = 1.036585507377728
$arrBF[14,3] = 1.03758035142347
$arrBF[14,4] = 1.041797662271323
$arrBF[15,0] = 1.02
$arrBF[15,1] = 1.034169375573884
$arrBF[15,2] = 1.036899797431976
$arrBF[15,3] = 1.037982009138677
$arrBF[15,4] = 1.042304606176919
$arrBF[16,0] = 1.02
$arrBF[16,1] = 1.034415921118967
$arrBF[16,2] = 1.037083106273685
$arrBF[16,3] = 1.038216319525841
$arrBF[16,4] = 1.042600324716593
$arrBF[17,0] = 1.02
$arrBF[17,1] = 1.03449998468072
$arrBF[17,2] = 1.037145607950174
$arrBF[17,3] = 1.038296218459368
$arrBF[17,4] = 1.042701161698508
$arrBF[18,0] = 1.02
$arrBF[18,1] = 1.034124024370717
$arrBF[18,2] = 1.036866078245775
$arrBF[18,3] = 1.037938911943576
$arrBF[18,4] = 1.042250213143681
$arrBF[19,0] = 1.02
$arrBF[19,1] = 1.032901590586333
$arrBF[19,2] = 1.035957167918415
$arrBF[19,3] = 1.036777627979934
$arrBF[19,4] = 1.040784445071151
$arrBF[20,0] = 1.02
$arrBF[20,1] = 1.032132876571266
$arrBF[20,2] = 1.035385597985415
$arrBF[20,3] = 1.036047751018824
$arrBF[20,4] = 1.039863084206879
$arrBF[21,0] = 1.02
$arrBF[21,1] = 1.032540398311486
$arrBF[21,2] = 1.035688607721678
$arrBF[21,3] = 1.03643464736524
$arrBF[21,4] = 1.040351494376055
$arrBF[22,0] = 1.02
$arrBF[22,1] = 1.034144516665811
$arrBF[22,2] = 1.03688131452915
$arrBF[22,3] = 1.037958385617272
$arrBF[22,4] = 1.042274790930889
$arrBF[23,0] = 1.02
$arrBF[23,1] = 1.036005517389412
$arrBF[23,2] = 1.038264944153457
$arrBF[23,3] = 1.039727787173971
$arrBF[23,4] = 1.044507726078981

$arrIN = New-Object "object[,]" 24,6
$arrIN[0,0] = 1.03897942308265
$arrIN[0,1] = 1.042592086551811
$arrIN[0,2] = 1.042153558752941
$arrIN[0,3] = 1.043920395528365
$arrIN[0,4] = 1.049055493391111
$arrIN[0,5] = 1.018061510580346
$arrIN[1,0] = 1.039286177178536
$arrIN[1,1] = 1.04331316214976
$arrIN[1,2] = 1.04276462806853
$arrIN[1,3] = 1.044756394764788
$arrIN[1,4] = 1.05016026266764
$arrIN[1,5] = 1.018307140185712
$arrIN[2,0] = 1.039483451324692
$arrIN[2,1] = 1.043779109704424
$arrIN[2,2] = 1.043159286619796
$arrIN[2,3] = 1.045297010382684
$arrIN[2,4] = 1.050874858251994
$arrIN[2,5] = 1.018465712321608
$arrIN[3,0] = 1.039566094315913
$arrIN[3,1] = 1.043974841953926
$arrIN[3,2] = 1.043325023043193
$arrIN[3,3] = 1.045524206046213
$arrIN[3,4] = 1.051175212391989
$arrIN[3,5] = 1.018532288256329
$arrIN[4,0] = 1.039579953373315
$arrIN[4,1] = 1.044007697349101
$arrIN[4,2] = 1.043352840481625
$arrIN[4,3] = 1.045562348565306
$arrIN[4,4] = 1.051225639586815
$arrIN[4,5] = 1.018543461500125
$arrIN[5,0] = 1.039484556747902
$arrIN[5,1] = 1.043781725685608
$arrIN[5,2] = 1.043161501897387
$arrIN[5,3] = 1.045300046492868
$arrIN[5,4] = 1.050878871840005
$arrIN[5,5] = 1.018466602257484
$arrIN[6,0] = 1.039083343971687
$arrIN[6,1] = 1.042835909410763
$arrIN[6,2] = 1.042360226645169
$arrIN[6,3] = 1.044202994257957
$arrIN[6,4] = 1.049428910824578
$arrIN[6,5] = 1.018144598162248
$arrIN[7,0] = 1.038367035207969
$arrIN[7,1] = 1.041164371285911
$arrIN[7,2] = 1.040942572062459
$arrIN[7,3] = 1.042267293544369
$arrIN[7,4] = 1.046871811134638
$arrIN[7,5] = 1.017574378239314
$arrIN[8,0] = 1.037883220448114
$arrIN[8,1] = 1.040046704099692
$arrIN[8,2] = 1.039993622737521
$arrIN[8,3] = 1.040975087591861
$arrIN[8,4] = 1.045165589783451
$arrIN[8,5] = 1.017192343629473
$arrIN[9,0] = 1.037672232406196
$arrIN[9,1] = 1.039561951665412
$arrIN[9,2] = 1.039581801899764
$arrIN[9,3] = 1.040415130172759
$arrIN[9,4] = 1.044426403906922
$arrIN[9,5] = 1.017026470216879
$arrIN[10,0] = 1.037593637464846
$arrIN[10,1] = 1.03938177299703
$arrIN[10,2] = 1.039428694808214
$arrIN[10,3] = 1.040207072819945
$arrIN[10,4] = 1.044151778317237
$arrIN[10,5] = 1.016964789777993
$arrIN[11,0] = 1.037610506514042
$arrIN[11,1] = 1.03942042737377
$arrIN[11,2] = 1.039461543071945
$arrIN[11,3] = 1.040251704744453
$arrIN[11,4] = 1.044210689139768
$arrIN[11,5] = 1.016978023504843
$arrIN[12,0] = 1.037665740306938
$arrIN[12,1] = 1.03954706049234
$arrIN[12,2] = 1.039569148847114
$arrIN[12,3] = 1.040397933400791
$arrIN[12,4] = 1.044403704485084
$arrIN[12,5] = 1.0170213730734
$arrIN[13,0] = 1.037699741901392
$arrIN[13,1] = 1.039625067413944
$arrIN[13,2] = 1.039635429954104
$arrIN[13,3] = 1.040488021173046
$arrIN[13,4] = 1.044522619730965
$arrIN[13,5] = 1.017048073204394
$arrIN[14,0] = 1.037897191549231
$arrIN[14,1] = 1.040078858720733
$arrIN[14,2] = 1.040020934527287
$arrIN[14,3] = 1.0410122411452
$arrIN[14,4] = 1.045214638936472
$arrIN[14,5] = 1.017203342616233
$arrIN[15,0] = 1.038020646569819
$arrIN[15,1] = 1.040363296730856
$arrIN[15,2] = 1.040262504934452
$arrIN[15,3] = 1.041340956637084
$arrIN[15,4] = 1.045648621062077
$arrIN[15,5] = 1.017300618511344
$arrIN[16,0] = 1.038092511754477
$arrIN[16,1] = 1.040529127925786
$arrIN[16,2] = 1.040403320169255
$arrIN[16,3] = 1.04153264997858
$arrIN[16,4] = 1.045901718877752
$arrIN[16,5] = 1.017357314490917
$arrIN[17,0] = 1.038116991501889
$arrIN[17,1] = 1.040585659088034
$arrIN[17,2] = 1.040451319491426
$arrIN[17,3] = 1.041598005543987
$arrIN[17,4] = 1.045988012537435
$arrIN[17,5] = 1.017376638994931
$arrIN[18,0] = 1.038007415915127
$arrIN[18,1] = 1.040332787160006
$arrIN[18,2] = 1.040236595892678
$arrIN[18,3] = 1.041305692804392
$arrIN[18,4] = 1.045602062687523
$arrIN[18,5] = 1.01729018621841
$arrIN[19,0] = 1.037649481533526
$arrIN[19,1] = 1.039509773513153
$arrIN[19,2] = 1.039537465454309
$arrIN[19,3] = 1.040354874484001
$arrIN[19,4] = 1.044346867921291
$arrIN[19,5] = 1.017008609574492
$arrIN[20,0] = 1.037423134447667
$arrIN[20,1] = 1.038991617426631
$arrIN[20,2] = 1.039097092929956
$arrIN[20,3] = 1.039756685336145
$arrIN[20,4] = 1.043557334640014
$arrIN[20,5] = 1.016831179619102
$arrIN[21,0] = 1.037543248592454
$arrIN[21,1] = 1.03926636781471
$arrIN[21,2] = 1.039330618824218
$arrIN[21,3] = 1.040073832169495
$arrIN[21,4] = 1.043975914250591
$arrIN[21,5] = 1.016925275734629
$arrIN[22,0] = 1.038013394723416
$arrIN[22,1] = 1.040346573358524
$arrIN[22,2] = 1.040248303347302
$arrIN[22,3] = 1.041321627138438
$arrIN[22,4] = 1.045623100525536
$arrIN[22,5] = 1.017294900256589
$arrIN[23,0] = 1.038553323652858
$arrIN[23,1] = 1.041597085522584
$arrIN[23,2] = 1.041309746915057
$arrIN[23,3] = 1.042768023487114
$arrIN[23,4] = 1.047533138585178
$arrIN[23,5] = 1.017722126499095

$ws.Range("B2:F25").Value = $arrBF
$ws.Range("I2:N25").Value = $arrIN
Write-Host "Applied changes to rows 2-25"
